$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained 3 new data rows for "Sandia" (Hortaliza) at Feria Lagunitas de
# Puerto Montt, for a new price date (2021-12-16 / serial 44546), grouped as
# Primera / Segunda / Tercera quality - Región de O'Higgins origin.
# In the OOXML this shows up as the old rows 79-151 shifting down to 82-154
# (dimension A1:R151 -> A1:R154), which is exactly what inserting 3 rows above
# row 79 produces.
$ws.Rows("79:81").Insert()

$newRows = @(
    @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44546, 10, 100112028, "Sandia", "Sin especificar", "Primera", 300, 3500, 3500, 3500, "`$/unidad", "Región de O'Higgins", 3500, 1, "Hortaliza"),
    @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44546, 10, 100112028, "Sandia", "Sin especificar", "Segunda", 300, 2800, 2800, 2800, "`$/unidad", "Región de O'Higgins", 2800, 1, "Hortaliza"),
    @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44546, 10, 100112028, "Sandia", "Sin especificar", "Tercera", 300, 2500, 2500, 2500, "`$/unidad", "Región de O'Higgins", 2500, 1, "Hortaliza")
)

$startRow = 79
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowVals = $newRows[$r]
    $rowIndex = $startRow + $r
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowVals[$c]
    }
}
